$xmlFrag = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rStyle w:val="a9"/></w:rPr><w:t>hello</w:t></w:r></w:p>
'@
Write-Output $xmlFrag.Length
Write-Output $xmlFrag
